# Scheduled market-data refresh: recompute currentAveragePrice / Leve
# price / Leve profit columns (H:N) for the rows whose underlying
# Universalis market data changed since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2055.9473
$ws.Range("I86").Value = 1383
$ws.Range("K86").Value = 1383
$ws.Range("M86").Value = -260
$ws.Range("H89").Value = 2055.9473
$ws.Range("I89").Value = 1383
$ws.Range("K89").Value = 6915
$ws.Range("M89").Value = -1299
$ws.Range("H116").Value = 6171.6665
$ws.Range("I116").Value = 4995
$ws.Range("J116").Value = 6407
$ws.Range("K116").Value = 4995
$ws.Range("L116").Value = 6407
$ws.Range("M116").Value = -1553
$ws.Range("N116").Value = -13291
$ws.Range("H129").Value = 25000834
$ws.Range("I129").Value = 26316510
$ws.Range("K129").Value = 78949530
$ws.Range("M129").Value = -78944530
$ws.Range("H132").Value = 2160.276
$ws.Range("I132").Value = 2166
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6498
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3968
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 199776.5
$ws.Range("J136").Value = 199776.5
$ws.Range("L136").Value = 199776.5
$ws.Range("N136").Value = -209976.5
$ws.Range("H138").Value = 3248.0315
$ws.Range("I138").Value = 1861.1333
$ws.Range("J138").Value = 3508.075
$ws.Range("K138").Value = 5583.3999
$ws.Range("L138").Value = 10524.225
$ws.Range("M138").Value = -443.3999000000003
$ws.Range("N138").Value = -20804.225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4431.375
$ws.Range("I32").Value = 2887.8547
$ws.Range("K32").Value = 2887.8547
$ws.Range("M32").Value = -2600.8547
$ws.Range("H41").Value = 5664.6924
$ws.Range("I41").Value = 4894.8184
$ws.Range("J41").Value = 9899
$ws.Range("K41").Value = 4894.8184
$ws.Range("L41").Value = 9899
$ws.Range("M41").Value = -4480.8184
$ws.Range("N41").Value = -10727
$ws.Range("H61").Value = 1653.4
$ws.Range("I61").Value = 1653.4
$ws.Range("K61").Value = 1653.4
$ws.Range("M61").Value = -1441.4
$ws.Range("H74").Value = 64073.816
$ws.Range("I74").Value = 40485
$ws.Range("K74").Value = 40485
$ws.Range("M74").Value = -39611
$ws.Range("H77").Value = 64073.816
$ws.Range("I77").Value = 40485
$ws.Range("K77").Value = 202425
$ws.Range("M77").Value = -198057
$ws.Range("H102").Value = 4661.485
$ws.Range("I102").Value = 4289.5713
$ws.Range("K102").Value = 4289.5713
$ws.Range("M102").Value = -2667.5713
$ws.Range("H136").Value = 1653.4
$ws.Range("I136").Value = 1653.4
$ws.Range("K136").Value = 4960.200000000001
$ws.Range("M136").Value = -2410.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7089.884
$ws.Range("I86").Value = 5800.1514
$ws.Range("J86").Value = 11346
$ws.Range("K86").Value = 5800.1514
$ws.Range("L86").Value = 11346
$ws.Range("M86").Value = -4677.1514
$ws.Range("N86").Value = -13592
$ws.Range("H89").Value = 7089.884
$ws.Range("I89").Value = 5800.1514
$ws.Range("J89").Value = 11346
$ws.Range("K89").Value = 29000.757
$ws.Range("L89").Value = 56730
$ws.Range("M89").Value = -23384.757
$ws.Range("N89").Value = -67962
$ws.Range("H94").Value = 4220.5625
$ws.Range("I94").Value = 562.82355
$ws.Range("K94").Value = 562.82355
$ws.Range("M94").Value = -111.82355
$ws.Range("H104").Value = 28250
$ws.Range("J104").Value = 28250
$ws.Range("L104").Value = 28250
$ws.Range("N104").Value = -35238

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 160.54546
$ws.Range("J7").Value = 350
$ws.Range("L7").Value = 350
$ws.Range("N7").Value = -576
$ws.Range("H50").Value = 2500
$ws.Range("J50").Value = 2500
$ws.Range("L50").Value = 2500
$ws.Range("N50").Value = -3750
$ws.Range("H51").Value = 29998.5
$ws.Range("J51").Value = 29998.5
$ws.Range("L51").Value = 29998.5
$ws.Range("N51").Value = -31470.5
$ws.Range("H61").Value = 29998.5
$ws.Range("J61").Value = 29998.5
$ws.Range("L61").Value = 29998.5
$ws.Range("N61").Value = -30694.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 2810.7778
$ws.Range("I43").Value = 4499
$ws.Range("J43").Value = 2599.75
$ws.Range("K43").Value = 13497
$ws.Range("L43").Value = 7799.25
$ws.Range("M43").Value = -13383
$ws.Range("N43").Value = -8027.25
$ws.Range("H75").Value = 1010
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1010
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 3030
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -5026
$ws.Range("H78").Value = 1010
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1010
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 9090
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -19074
$ws.Range("H117").Value = 1828.2858
$ws.Range("J117").Value = 1809.8
$ws.Range("L117").Value = 5429.4
$ws.Range("N117").Value = -12313.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6832.222
$ws.Range("J70").Value = 7498.3335
$ws.Range("L70").Value = 7498.3335
$ws.Range("N70").Value = -8038.3335
$ws.Range("H73").Value = 6832.222
$ws.Range("J73").Value = 7498.3335
$ws.Range("L73").Value = 7498.3335
$ws.Range("N73").Value = -9370.333500000001
$ws.Range("H97").Value = 737.75757
$ws.Range("J97").Value = 721.4
$ws.Range("L97").Value = 721.4
$ws.Range("N97").Value = -1713.4
$ws.Range("H105").Value = 17666.666
$ws.Range("J105").Value = 17666.666
$ws.Range("L105").Value = 17666.666
$ws.Range("N105").Value = -24654.666
$ws.Range("H126").Value = 3262.5881
$ws.Range("I126").Value = 3074.6155
$ws.Range("J126").Value = 3873.5
$ws.Range("K126").Value = 9223.8465
$ws.Range("L126").Value = 11620.5
$ws.Range("M126").Value = -6753.8465
$ws.Range("N126").Value = -16560.5
$ws.Range("H132").Value = 5085.067
$ws.Range("I132").Value = 3227.7
$ws.Range("J132").Value = 8799.799999999999
$ws.Range("K132").Value = 9683.099999999999
$ws.Range("L132").Value = 26399.4
$ws.Range("M132").Value = -7153.099999999999
$ws.Range("N132").Value = -31459.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1588.75
$ws.Range("I55").Value = 1732.5
$ws.Range("J55").Value = 1349.1666
$ws.Range("K55").Value = 1732.5
$ws.Range("L55").Value = 1349.1666
$ws.Range("M55").Value = -1559.5
$ws.Range("N55").Value = -1695.1666
$ws.Range("H122").Value = 4688.1934
$ws.Range("I122").Value = 3529.5789
$ws.Range("J122").Value = 6522.6665
$ws.Range("K122").Value = 10588.7367
$ws.Range("L122").Value = 19567.9995
$ws.Range("M122").Value = -8138.736699999999
$ws.Range("N122").Value = -24467.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6961.5093
$ws.Range("I62").Value = 3291.1177
$ws.Range("K62").Value = 3291.1177
$ws.Range("M62").Value = -2667.1177
$ws.Range("H65").Value = 6961.5093
$ws.Range("I65").Value = 3291.1177
$ws.Range("K65").Value = 16455.5885
$ws.Range("M65").Value = -13335.5885
$ws.Range("H81").Value = 824.1111
$ws.Range("I81").Value = 859.8570999999999
$ws.Range("K81").Value = 1719.7142
$ws.Range("M81").Value = -658.7141999999999
$ws.Range("H84").Value = 824.1111
$ws.Range("I84").Value = 859.8570999999999
$ws.Range("K84").Value = 8598.571
$ws.Range("M84").Value = -3294.571
$ws.Range("H106").Value = 53799.5
$ws.Range("J106").Value = 53799.5
$ws.Range("L106").Value = 53799.5
$ws.Range("N106").Value = -56323.5
$ws.Range("H122").Value = 2989.9048
$ws.Range("I122").Value = 2261.4119
$ws.Range("K122").Value = 6784.2357
$ws.Range("M122").Value = -4334.2357
$ws.Range("H136").Value = 2056.838
$ws.Range("I136").Value = 1066.4231
$ws.Range("J136").Value = 4397.8184
$ws.Range("K136").Value = 3199.2693
$ws.Range("L136").Value = 13193.4552
$ws.Range("M136").Value = -649.2692999999999
$ws.Range("N136").Value = -18293.4552

